$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.288370847702026
$ws.Range("B1").Value = 2.316179513931274
$ws.Range("D1").Value = 1.377618312835693
$ws.Range("E1").Value = 0.8362594842910767
